$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148, shifting existing rows 148-200 down to 149-201.
$ws.Rows(148).Insert()

# Populate the newly inserted row 148 with the new record's data.
$ws.Cells.Item(148, 1).Value = 5
$ws.Cells.Item(148, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(148, 3).Value = "Maule"
$ws.Cells.Item(148, 4).Value = (Get-Date -Year 2022 -Month 3 -Day 7).Date
$ws.Cells.Item(148, 5).Value = 7
$ws.Cells.Item(148, 6).Value = 100112008
$ws.Cells.Item(148, 7).Value = "Coliflor"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 2000
$ws.Cells.Item(148, 11).Value = 1600
$ws.Cells.Item(148, 12).Value = 1600
$ws.Cells.Item(148, 13).Value = 1600
$ws.Cells.Item(148, 14).Value = "`$/unidad"
$ws.Cells.Item(148, 15).Value = "Región del Maule"
$ws.Cells.Item(148, 16).Value = 1600
$ws.Cells.Item(148, 17).Value = 1
$ws.Cells.Item(148, 18).Value = "Hortaliza"
